$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-assert pre-existing blank numeric cells (rows 1-27) as true blanks.
#     The engine round-trips an untouched blank <v></v> numeric cell as 0 on
#     save; explicitly nulling preserves blank semantics through the save. ---
$ws.Cells.Item(2,8).Value = $null
$ws.Cells.Item(2,9).Value = $null
$ws.Cells.Item(2,10).Value = $null
$ws.Cells.Item(2,11).Value = $null
$ws.Cells.Item(2,12).Value = $null
$ws.Cells.Item(2,13).Value = $null
$ws.Cells.Item(2,14).Value = $null
$ws.Cells.Item(2,16).Value = $null
$ws.Cells.Item(2,17).Value = $null
$ws.Cells.Item(2,18).Value = $null
$ws.Cells.Item(2,28).Value = $null
$ws.Cells.Item(3,8).Value = $null
$ws.Cells.Item(3,9).Value = $null
$ws.Cells.Item(3,10).Value = $null
$ws.Cells.Item(3,11).Value = $null
$ws.Cells.Item(3,12).Value = $null
$ws.Cells.Item(3,13).Value = $null
$ws.Cells.Item(3,14).Value = $null
$ws.Cells.Item(3,16).Value = $null
$ws.Cells.Item(3,17).Value = $null
$ws.Cells.Item(3,18).Value = $null
$ws.Cells.Item(3,28).Value = $null
$ws.Cells.Item(4,8).Value = $null
$ws.Cells.Item(4,9).Value = $null
$ws.Cells.Item(4,10).Value = $null
$ws.Cells.Item(4,11).Value = $null
$ws.Cells.Item(4,12).Value = $null
$ws.Cells.Item(4,13).Value = $null
$ws.Cells.Item(4,14).Value = $null
$ws.Cells.Item(4,16).Value = $null
$ws.Cells.Item(4,17).Value = $null
$ws.Cells.Item(4,18).Value = $null
$ws.Cells.Item(4,28).Value = $null
$ws.Cells.Item(5,8).Value = $null
$ws.Cells.Item(5,9).Value = $null
$ws.Cells.Item(5,10).Value = $null
$ws.Cells.Item(5,11).Value = $null
$ws.Cells.Item(5,12).Value = $null
$ws.Cells.Item(5,13).Value = $null
$ws.Cells.Item(5,14).Value = $null
$ws.Cells.Item(5,16).Value = $null
$ws.Cells.Item(5,17).Value = $null
$ws.Cells.Item(5,18).Value = $null
$ws.Cells.Item(5,28).Value = $null
$ws.Cells.Item(6,17).Value = $null
$ws.Cells.Item(6,18).Value = $null
$ws.Cells.Item(6,28).Value = $null
$ws.Cells.Item(7,8).Value = $null
$ws.Cells.Item(7,11).Value = $null
$ws.Cells.Item(7,12).Value = $null
$ws.Cells.Item(7,17).Value = $null
$ws.Cells.Item(7,29).Value = $null
$ws.Cells.Item(7,30).Value = $null
$ws.Cells.Item(7,31).Value = $null
$ws.Cells.Item(8,8).Value = $null
$ws.Cells.Item(8,9).Value = $null
$ws.Cells.Item(8,10).Value = $null
$ws.Cells.Item(8,11).Value = $null
$ws.Cells.Item(8,12).Value = $null
$ws.Cells.Item(8,13).Value = $null
$ws.Cells.Item(8,14).Value = $null
$ws.Cells.Item(9,17).Value = $null
$ws.Cells.Item(9,18).Value = $null
$ws.Cells.Item(10,8).Value = $null
$ws.Cells.Item(10,11).Value = $null
$ws.Cells.Item(10,12).Value = $null
$ws.Cells.Item(10,17).Value = $null
$ws.Cells.Item(10,29).Value = $null
$ws.Cells.Item(10,30).Value = $null
$ws.Cells.Item(10,31).Value = $null
$ws.Cells.Item(11,8).Value = $null
$ws.Cells.Item(11,9).Value = $null
$ws.Cells.Item(11,10).Value = $null
$ws.Cells.Item(11,11).Value = $null
$ws.Cells.Item(11,12).Value = $null
$ws.Cells.Item(11,13).Value = $null
$ws.Cells.Item(11,14).Value = $null
$ws.Cells.Item(12,8).Value = $null
$ws.Cells.Item(12,9).Value = $null
$ws.Cells.Item(12,10).Value = $null
$ws.Cells.Item(12,11).Value = $null
$ws.Cells.Item(12,12).Value = $null
$ws.Cells.Item(12,13).Value = $null
$ws.Cells.Item(12,14).Value = $null
$ws.Cells.Item(12,17).Value = $null
$ws.Cells.Item(12,18).Value = $null
$ws.Cells.Item(13,8).Value = $null
$ws.Cells.Item(13,9).Value = $null
$ws.Cells.Item(13,10).Value = $null
$ws.Cells.Item(13,11).Value = $null
$ws.Cells.Item(13,12).Value = $null
$ws.Cells.Item(13,13).Value = $null
$ws.Cells.Item(13,14).Value = $null
$ws.Cells.Item(13,17).Value = $null
$ws.Cells.Item(13,18).Value = $null
$ws.Cells.Item(14,8).Value = $null
$ws.Cells.Item(14,9).Value = $null
$ws.Cells.Item(14,10).Value = $null
$ws.Cells.Item(14,11).Value = $null
$ws.Cells.Item(14,12).Value = $null
$ws.Cells.Item(14,13).Value = $null
$ws.Cells.Item(14,14).Value = $null
$ws.Cells.Item(14,17).Value = $null
$ws.Cells.Item(14,18).Value = $null
$ws.Cells.Item(15,17).Value = $null
$ws.Cells.Item(15,18).Value = $null
$ws.Cells.Item(16,8).Value = $null
$ws.Cells.Item(16,11).Value = $null
$ws.Cells.Item(16,12).Value = $null
$ws.Cells.Item(16,17).Value = $null
$ws.Cells.Item(16,29).Value = $null
$ws.Cells.Item(16,30).Value = $null
$ws.Cells.Item(16,31).Value = $null
$ws.Cells.Item(17,8).Value = $null
$ws.Cells.Item(17,9).Value = $null
$ws.Cells.Item(17,10).Value = $null
$ws.Cells.Item(17,11).Value = $null
$ws.Cells.Item(17,12).Value = $null
$ws.Cells.Item(17,13).Value = $null
$ws.Cells.Item(17,14).Value = $null
$ws.Cells.Item(18,17).Value = $null
$ws.Cells.Item(18,18).Value = $null
$ws.Cells.Item(19,17).Value = $null
$ws.Cells.Item(19,18).Value = $null
$ws.Cells.Item(20,8).Value = $null
$ws.Cells.Item(20,11).Value = $null
$ws.Cells.Item(20,12).Value = $null
$ws.Cells.Item(20,17).Value = $null
$ws.Cells.Item(20,29).Value = $null
$ws.Cells.Item(20,30).Value = $null
$ws.Cells.Item(20,31).Value = $null
$ws.Cells.Item(21,8).Value = $null
$ws.Cells.Item(21,9).Value = $null
$ws.Cells.Item(21,10).Value = $null
$ws.Cells.Item(21,11).Value = $null
$ws.Cells.Item(21,12).Value = $null
$ws.Cells.Item(21,13).Value = $null
$ws.Cells.Item(21,14).Value = $null
$ws.Cells.Item(22,17).Value = $null
$ws.Cells.Item(22,18).Value = $null
$ws.Cells.Item(23,8).Value = $null
$ws.Cells.Item(23,11).Value = $null
$ws.Cells.Item(23,12).Value = $null
$ws.Cells.Item(23,17).Value = $null
$ws.Cells.Item(23,29).Value = $null
$ws.Cells.Item(23,30).Value = $null
$ws.Cells.Item(23,31).Value = $null
$ws.Cells.Item(24,8).Value = $null
$ws.Cells.Item(24,9).Value = $null
$ws.Cells.Item(24,10).Value = $null
$ws.Cells.Item(24,11).Value = $null
$ws.Cells.Item(24,12).Value = $null
$ws.Cells.Item(24,13).Value = $null
$ws.Cells.Item(24,14).Value = $null
$ws.Cells.Item(25,17).Value = $null
$ws.Cells.Item(25,18).Value = $null
$ws.Cells.Item(26,8).Value = $null
$ws.Cells.Item(26,11).Value = $null
$ws.Cells.Item(26,12).Value = $null
$ws.Cells.Item(26,17).Value = $null
$ws.Cells.Item(26,29).Value = $null
$ws.Cells.Item(26,30).Value = $null
$ws.Cells.Item(26,31).Value = $null
$ws.Cells.Item(27,8).Value = $null
$ws.Cells.Item(27,9).Value = $null
$ws.Cells.Item(27,10).Value = $null
$ws.Cells.Item(27,11).Value = $null
$ws.Cells.Item(27,12).Value = $null
$ws.Cells.Item(27,13).Value = $null
$ws.Cells.Item(27,14).Value = $null

# --- Modify existing rows 25-27: high/close updates, row 27 rebound pct ---
$ws.Cells.Item(25,3).Value = 1.243
$ws.Cells.Item(25,5).Value = 1.2
$ws.Cells.Item(26,3).Value = 1.243
$ws.Cells.Item(26,5).Value = 1.2
$ws.Cells.Item(27,3).Value = 1.243
$ws.Cells.Item(27,5).Value = 1.2
$ws.Cells.Item(27,17).Value = 21.862745

# --- Append new rows 28-41 ---

# Row 28
$ws.Cells.Item(28,1).NumberFormat = "@"
$ws.Cells.Item(28,1).Value = '2025-10-19'
$ws.Cells.Item(28,1).ClearFormats()
$ws.Cells.Item(28,2).Value = 1.2
$ws.Cells.Item(28,3).Value = 1.26
$ws.Cells.Item(28,4).Value = 1.12
$ws.Cells.Item(28,5).Value = 1.158
$ws.Cells.Item(28,6).Value = 'wait'
$ws.Cells.Item(28,7).Value = $true
$ws.Cells.Item(28,8).Value = 2
$ws.Cells.Item(28,9).Value = 'BUY B2'
$ws.Cells.Item(28,10).Value = 'LOW'
$ws.Cells.Item(28,11).Value = 'B2'
$ws.Cells.Item(28,12).Value = 1.1388
$ws.Cells.Item(28,13).Value = 1.12
$ws.Cells.Item(28,14).Value = 1.1388
$ws.Cells.Item(28,15).Value = 2.19
$ws.Cells.Item(28,16).Value = 1.02
$ws.Cells.Item(28,17).Value = $null
$ws.Cells.Item(28,18).Value = $null
$ws.Cells.Item(28,19).Value = 6
$ws.Cells.Item(28,20).Value = 1.2264
$ws.Cells.Item(28,21).Value = 1.1388
$ws.Cells.Item(28,22).Value = 1.0074
$ws.Cells.Item(28,23).Value = 0.8979
$ws.Cells.Item(28,24).Value = 0.7665
$ws.Cells.Item(28,25).Value = 0.6132
$ws.Cells.Item(28,26).Value = 0.4599
$ws.Cells.Item(28,27).Value = 0.4161
$ws.Cells.Item(28,28).Value = 1.19646
$ws.Cells.Item(28,29).Value = 'B2'
$ws.Cells.Item(28,30).Value = 1.1388
$ws.Cells.Item(28,31).Value = 1.12

# Row 29
$ws.Cells.Item(29,1).NumberFormat = "@"
$ws.Cells.Item(29,1).Value = '2025-10-19'
$ws.Cells.Item(29,1).ClearFormats()
$ws.Cells.Item(29,2).Value = 1.2
$ws.Cells.Item(29,3).Value = 1.26
$ws.Cells.Item(29,4).Value = 1.12
$ws.Cells.Item(29,5).Value = 1.158
$ws.Cells.Item(29,6).Value = 'wait'
$ws.Cells.Item(29,7).Value = $false
$ws.Cells.Item(29,8).Value = $null
$ws.Cells.Item(29,9).Value = 'SELL S'
$ws.Cells.Item(29,10).Value = 'HIGH'
$ws.Cells.Item(29,11).Value = $null
$ws.Cells.Item(29,12).Value = $null
$ws.Cells.Item(29,13).Value = 1.19646
$ws.Cells.Item(29,14).Value = 1.19646
$ws.Cells.Item(29,15).Value = 2.19
$ws.Cells.Item(29,16).Value = 1.02
$ws.Cells.Item(29,17).Value = $null
$ws.Cells.Item(29,18).Value = 17.3
$ws.Cells.Item(29,19).Value = 6
$ws.Cells.Item(29,20).Value = 1.2264
$ws.Cells.Item(29,21).Value = 1.1388
$ws.Cells.Item(29,22).Value = 1.0074
$ws.Cells.Item(29,23).Value = 0.8979
$ws.Cells.Item(29,24).Value = 0.7665
$ws.Cells.Item(29,25).Value = 0.6132
$ws.Cells.Item(29,26).Value = 0.4599
$ws.Cells.Item(29,27).Value = 0.4161
$ws.Cells.Item(29,28).Value = 1.19646
$ws.Cells.Item(29,29).Value = $null
$ws.Cells.Item(29,30).Value = $null
$ws.Cells.Item(29,31).Value = $null

# Row 30
$ws.Cells.Item(30,1).NumberFormat = "@"
$ws.Cells.Item(30,1).Value = '2025-10-19'
$ws.Cells.Item(30,1).ClearFormats()
$ws.Cells.Item(30,2).Value = 1.2
$ws.Cells.Item(30,3).Value = 1.26
$ws.Cells.Item(30,4).Value = 1.12
$ws.Cells.Item(30,5).Value = 1.158
$ws.Cells.Item(30,6).Value = 'wait'
$ws.Cells.Item(30,7).Value = $false
$ws.Cells.Item(30,8).Value = $null
$ws.Cells.Item(30,9).Value = $null
$ws.Cells.Item(30,10).Value = $null
$ws.Cells.Item(30,11).Value = $null
$ws.Cells.Item(30,12).Value = $null
$ws.Cells.Item(30,13).Value = $null
$ws.Cells.Item(30,14).Value = $null
$ws.Cells.Item(30,15).Value = 2.19
$ws.Cells.Item(30,16).Value = 1.02
$ws.Cells.Item(30,17).Value = 23.529412
$ws.Cells.Item(30,18).Value = 17.3
$ws.Cells.Item(30,19).Value = 6
$ws.Cells.Item(30,20).Value = 1.2264
$ws.Cells.Item(30,21).Value = 1.1388
$ws.Cells.Item(30,22).Value = 1.0074
$ws.Cells.Item(30,23).Value = 0.8979
$ws.Cells.Item(30,24).Value = 0.7665
$ws.Cells.Item(30,25).Value = 0.6132
$ws.Cells.Item(30,26).Value = 0.4599
$ws.Cells.Item(30,27).Value = 0.4161
$ws.Cells.Item(30,28).Value = 1.19646
$ws.Cells.Item(30,29).Value = 'B2'
$ws.Cells.Item(30,30).Value = 1.1388
$ws.Cells.Item(30,31).Value = 1.12

# Row 31
$ws.Cells.Item(31,1).NumberFormat = "@"
$ws.Cells.Item(31,1).Value = '2025-10-20'
$ws.Cells.Item(31,1).ClearFormats()
$ws.Cells.Item(31,2).Value = 1.158
$ws.Cells.Item(31,3).Value = 1.192
$ws.Cells.Item(31,4).Value = 1.024
$ws.Cells.Item(31,5).Value = 1.042
$ws.Cells.Item(31,6).Value = 'wait'
$ws.Cells.Item(31,7).Value = $true
$ws.Cells.Item(31,8).Value = 2
$ws.Cells.Item(31,9).Value = 'BUY B2'
$ws.Cells.Item(31,10).Value = 'LOW'
$ws.Cells.Item(31,11).Value = 'B2'
$ws.Cells.Item(31,12).Value = 1.1388
$ws.Cells.Item(31,13).Value = 1.024
$ws.Cells.Item(31,14).Value = 1.1388
$ws.Cells.Item(31,15).Value = 2.19
$ws.Cells.Item(31,16).Value = 1.02
$ws.Cells.Item(31,17).Value = $null
$ws.Cells.Item(31,18).Value = $null
$ws.Cells.Item(31,19).Value = 6
$ws.Cells.Item(31,20).Value = 1.2264
$ws.Cells.Item(31,21).Value = 1.1388
$ws.Cells.Item(31,22).Value = 1.0074
$ws.Cells.Item(31,23).Value = 0.8979
$ws.Cells.Item(31,24).Value = 0.7665
$ws.Cells.Item(31,25).Value = 0.6132
$ws.Cells.Item(31,26).Value = 0.4599
$ws.Cells.Item(31,27).Value = 0.4161
$ws.Cells.Item(31,28).Value = 1.19646
$ws.Cells.Item(31,29).Value = 'B2'
$ws.Cells.Item(31,30).Value = 1.1388
$ws.Cells.Item(31,31).Value = 1.024

# Row 32
$ws.Cells.Item(32,1).NumberFormat = "@"
$ws.Cells.Item(32,1).Value = '2025-10-20'
$ws.Cells.Item(32,1).ClearFormats()
$ws.Cells.Item(32,2).Value = 1.158
$ws.Cells.Item(32,3).Value = 1.192
$ws.Cells.Item(32,4).Value = 1.024
$ws.Cells.Item(32,5).Value = 1.042
$ws.Cells.Item(32,6).Value = 'wait'
$ws.Cells.Item(32,7).Value = $true
$ws.Cells.Item(32,8).Value = 2
$ws.Cells.Item(32,9).Value = $null
$ws.Cells.Item(32,10).Value = $null
$ws.Cells.Item(32,11).Value = $null
$ws.Cells.Item(32,12).Value = $null
$ws.Cells.Item(32,13).Value = $null
$ws.Cells.Item(32,14).Value = $null
$ws.Cells.Item(32,15).Value = 2.19
$ws.Cells.Item(32,16).Value = 1.02
$ws.Cells.Item(32,17).Value = 16.862745
$ws.Cells.Item(32,18).Value = 17.3
$ws.Cells.Item(32,19).Value = 6
$ws.Cells.Item(32,20).Value = 1.2264
$ws.Cells.Item(32,21).Value = 1.1388
$ws.Cells.Item(32,22).Value = 1.0074
$ws.Cells.Item(32,23).Value = 0.8979
$ws.Cells.Item(32,24).Value = 0.7665
$ws.Cells.Item(32,25).Value = 0.6132
$ws.Cells.Item(32,26).Value = 0.4599
$ws.Cells.Item(32,27).Value = 0.4161
$ws.Cells.Item(32,28).Value = 1.19646
$ws.Cells.Item(32,29).Value = 'B2'
$ws.Cells.Item(32,30).Value = 1.1388
$ws.Cells.Item(32,31).Value = 1.024

# Row 33
$ws.Cells.Item(33,1).NumberFormat = "@"
$ws.Cells.Item(33,1).Value = '2025-10-21'
$ws.Cells.Item(33,1).ClearFormats()
$ws.Cells.Item(33,2).Value = 1.042
$ws.Cells.Item(33,3).Value = 1.063
$ws.Cells.Item(33,4).Value = 0.929
$ws.Cells.Item(33,5).Value = 0.954
$ws.Cells.Item(33,6).Value = 'wait'
$ws.Cells.Item(33,7).Value = $true
$ws.Cells.Item(33,8).Value = 3
$ws.Cells.Item(33,9).Value = 'ADD B3'
$ws.Cells.Item(33,10).Value = 'LOW'
$ws.Cells.Item(33,11).Value = 'B3'
$ws.Cells.Item(33,12).Value = 1.0074
$ws.Cells.Item(33,13).Value = 0.929
$ws.Cells.Item(33,14).Value = 1.0074
$ws.Cells.Item(33,15).Value = 2.19
$ws.Cells.Item(33,16).Value = 0.929
$ws.Cells.Item(33,17).Value = $null
$ws.Cells.Item(33,18).Value = $null
$ws.Cells.Item(33,19).Value = 6
$ws.Cells.Item(33,20).Value = 1.2264
$ws.Cells.Item(33,21).Value = 1.1388
$ws.Cells.Item(33,22).Value = 1.0074
$ws.Cells.Item(33,23).Value = 0.8979
$ws.Cells.Item(33,24).Value = 0.7665
$ws.Cells.Item(33,25).Value = 0.6132
$ws.Cells.Item(33,26).Value = 0.4599
$ws.Cells.Item(33,27).Value = 0.4161
$ws.Cells.Item(33,28).Value = 1.19646
$ws.Cells.Item(33,29).Value = 'B3'
$ws.Cells.Item(33,30).Value = 1.0074
$ws.Cells.Item(33,31).Value = 0.929

# Row 34
$ws.Cells.Item(34,1).NumberFormat = "@"
$ws.Cells.Item(34,1).Value = '2025-10-21'
$ws.Cells.Item(34,1).ClearFormats()
$ws.Cells.Item(34,2).Value = 1.042
$ws.Cells.Item(34,3).Value = 1.063
$ws.Cells.Item(34,4).Value = 0.929
$ws.Cells.Item(34,5).Value = 0.954
$ws.Cells.Item(34,6).Value = 'wait'
$ws.Cells.Item(34,7).Value = $true
$ws.Cells.Item(34,8).Value = 3
$ws.Cells.Item(34,9).Value = $null
$ws.Cells.Item(34,10).Value = $null
$ws.Cells.Item(34,11).Value = $null
$ws.Cells.Item(34,12).Value = $null
$ws.Cells.Item(34,13).Value = $null
$ws.Cells.Item(34,14).Value = $null
$ws.Cells.Item(34,15).Value = 2.19
$ws.Cells.Item(34,16).Value = 0.929
$ws.Cells.Item(34,17).Value = 14.424112
$ws.Cells.Item(34,18).Value = 24.4
$ws.Cells.Item(34,19).Value = 6
$ws.Cells.Item(34,20).Value = 1.2264
$ws.Cells.Item(34,21).Value = 1.1388
$ws.Cells.Item(34,22).Value = 1.0074
$ws.Cells.Item(34,23).Value = 0.8979
$ws.Cells.Item(34,24).Value = 0.7665
$ws.Cells.Item(34,25).Value = 0.6132
$ws.Cells.Item(34,26).Value = 0.4599
$ws.Cells.Item(34,27).Value = 0.4161
$ws.Cells.Item(34,28).Value = 1.19646
$ws.Cells.Item(34,29).Value = 'B3'
$ws.Cells.Item(34,30).Value = 1.0074
$ws.Cells.Item(34,31).Value = 0.929

# Row 35
$ws.Cells.Item(35,1).NumberFormat = "@"
$ws.Cells.Item(35,1).Value = '2025-10-22'
$ws.Cells.Item(35,1).ClearFormats()
$ws.Cells.Item(35,2).Value = 0.955
$ws.Cells.Item(35,3).Value = 1.118
$ws.Cells.Item(35,4).Value = 0.949
$ws.Cells.Item(35,5).Value = 1.106
$ws.Cells.Item(35,6).Value = 'wait'
$ws.Cells.Item(35,7).Value = $true
$ws.Cells.Item(35,8).Value = 3
$ws.Cells.Item(35,9).Value = $null
$ws.Cells.Item(35,10).Value = $null
$ws.Cells.Item(35,11).Value = $null
$ws.Cells.Item(35,12).Value = $null
$ws.Cells.Item(35,13).Value = $null
$ws.Cells.Item(35,14).Value = $null
$ws.Cells.Item(35,15).Value = 2.19
$ws.Cells.Item(35,16).Value = 0.929
$ws.Cells.Item(35,17).Value = 20.344456
$ws.Cells.Item(35,18).Value = 24.4
$ws.Cells.Item(35,19).Value = 6
$ws.Cells.Item(35,20).Value = 1.2264
$ws.Cells.Item(35,21).Value = 1.1388
$ws.Cells.Item(35,22).Value = 1.0074
$ws.Cells.Item(35,23).Value = 0.8979
$ws.Cells.Item(35,24).Value = 0.7665
$ws.Cells.Item(35,25).Value = 0.6132
$ws.Cells.Item(35,26).Value = 0.4599
$ws.Cells.Item(35,27).Value = 0.4161
$ws.Cells.Item(35,28).Value = 1.19646
$ws.Cells.Item(35,29).Value = 'B3'
$ws.Cells.Item(35,30).Value = 1.0074
$ws.Cells.Item(35,31).Value = 0.949

# Row 36
$ws.Cells.Item(36,1).NumberFormat = "@"
$ws.Cells.Item(36,1).Value = '2025-10-23'
$ws.Cells.Item(36,1).ClearFormats()
$ws.Cells.Item(36,2).Value = 1.105
$ws.Cells.Item(36,3).Value = 1.154
$ws.Cells.Item(36,4).Value = 1.048
$ws.Cells.Item(36,5).Value = 1.119
$ws.Cells.Item(36,6).Value = 'wait'
$ws.Cells.Item(36,7).Value = $true
$ws.Cells.Item(36,8).Value = 3
$ws.Cells.Item(36,9).Value = $null
$ws.Cells.Item(36,10).Value = $null
$ws.Cells.Item(36,11).Value = $null
$ws.Cells.Item(36,12).Value = $null
$ws.Cells.Item(36,13).Value = $null
$ws.Cells.Item(36,14).Value = $null
$ws.Cells.Item(36,15).Value = 2.19
$ws.Cells.Item(36,16).Value = 0.929
$ws.Cells.Item(36,17).Value = 24.219591
$ws.Cells.Item(36,18).Value = 24.4
$ws.Cells.Item(36,19).Value = 6
$ws.Cells.Item(36,20).Value = 1.2264
$ws.Cells.Item(36,21).Value = 1.1388
$ws.Cells.Item(36,22).Value = 1.0074
$ws.Cells.Item(36,23).Value = 0.8979
$ws.Cells.Item(36,24).Value = 0.7665
$ws.Cells.Item(36,25).Value = 0.6132
$ws.Cells.Item(36,26).Value = 0.4599
$ws.Cells.Item(36,27).Value = 0.4161
$ws.Cells.Item(36,28).Value = 1.19646
$ws.Cells.Item(36,29).Value = 'B2'
$ws.Cells.Item(36,30).Value = 1.1388
$ws.Cells.Item(36,31).Value = 1.048

# Row 37
$ws.Cells.Item(37,1).NumberFormat = "@"
$ws.Cells.Item(37,1).Value = '2025-10-24'
$ws.Cells.Item(37,1).ClearFormats()
$ws.Cells.Item(37,2).Value = 1.118
$ws.Cells.Item(37,3).Value = 1.161
$ws.Cells.Item(37,4).Value = 1.091
$ws.Cells.Item(37,5).Value = 1.14
$ws.Cells.Item(37,6).Value = 'wait'
$ws.Cells.Item(37,7).Value = $false
$ws.Cells.Item(37,8).Value = $null
$ws.Cells.Item(37,9).Value = 'SELL S'
$ws.Cells.Item(37,10).Value = 'HIGH'
$ws.Cells.Item(37,11).Value = $null
$ws.Cells.Item(37,12).Value = $null
$ws.Cells.Item(37,13).Value = 1.155676
$ws.Cells.Item(37,14).Value = 1.155676
$ws.Cells.Item(37,15).Value = 2.19
$ws.Cells.Item(37,16).Value = 0.929
$ws.Cells.Item(37,17).Value = $null
$ws.Cells.Item(37,18).Value = 24.4
$ws.Cells.Item(37,19).Value = 6
$ws.Cells.Item(37,20).Value = 1.2264
$ws.Cells.Item(37,21).Value = 1.1388
$ws.Cells.Item(37,22).Value = 1.0074
$ws.Cells.Item(37,23).Value = 0.8979
$ws.Cells.Item(37,24).Value = 0.7665
$ws.Cells.Item(37,25).Value = 0.6132
$ws.Cells.Item(37,26).Value = 0.4599
$ws.Cells.Item(37,27).Value = 0.4161
$ws.Cells.Item(37,28).Value = 1.155676
$ws.Cells.Item(37,29).Value = $null
$ws.Cells.Item(37,30).Value = $null
$ws.Cells.Item(37,31).Value = $null

# Row 38
$ws.Cells.Item(38,1).NumberFormat = "@"
$ws.Cells.Item(38,1).Value = '2025-10-24'
$ws.Cells.Item(38,1).ClearFormats()
$ws.Cells.Item(38,2).Value = 1.118
$ws.Cells.Item(38,3).Value = 1.161
$ws.Cells.Item(38,4).Value = 1.091
$ws.Cells.Item(38,5).Value = 1.14
$ws.Cells.Item(38,6).Value = 'wait'
$ws.Cells.Item(38,7).Value = $false
$ws.Cells.Item(38,8).Value = $null
$ws.Cells.Item(38,9).Value = $null
$ws.Cells.Item(38,10).Value = $null
$ws.Cells.Item(38,11).Value = $null
$ws.Cells.Item(38,12).Value = $null
$ws.Cells.Item(38,13).Value = $null
$ws.Cells.Item(38,14).Value = $null
$ws.Cells.Item(38,15).Value = 2.19
$ws.Cells.Item(38,16).Value = 0.929
$ws.Cells.Item(38,17).Value = 24.973089
$ws.Cells.Item(38,18).Value = 24.4
$ws.Cells.Item(38,19).Value = 6
$ws.Cells.Item(38,20).Value = 1.2264
$ws.Cells.Item(38,21).Value = 1.1388
$ws.Cells.Item(38,22).Value = 1.0074
$ws.Cells.Item(38,23).Value = 0.8979
$ws.Cells.Item(38,24).Value = 0.7665
$ws.Cells.Item(38,25).Value = 0.6132
$ws.Cells.Item(38,26).Value = 0.4599
$ws.Cells.Item(38,27).Value = 0.4161
$ws.Cells.Item(38,28).Value = 1.155676
$ws.Cells.Item(38,29).Value = 'B2'
$ws.Cells.Item(38,30).Value = 1.1388
$ws.Cells.Item(38,31).Value = 1.091

# Row 39
$ws.Cells.Item(39,1).NumberFormat = "@"
$ws.Cells.Item(39,1).Value = '2025-10-25'
$ws.Cells.Item(39,1).ClearFormats()
$ws.Cells.Item(39,2).Value = 1.139
$ws.Cells.Item(39,3).Value = 1.22
$ws.Cells.Item(39,4).Value = 1.12
$ws.Cells.Item(39,5).Value = 1.131
$ws.Cells.Item(39,6).Value = 'wait'
$ws.Cells.Item(39,7).Value = $true
$ws.Cells.Item(39,8).Value = 2
$ws.Cells.Item(39,9).Value = 'BUY B2'
$ws.Cells.Item(39,10).Value = 'LOW'
$ws.Cells.Item(39,11).Value = 'B2'
$ws.Cells.Item(39,12).Value = 1.1388
$ws.Cells.Item(39,13).Value = 1.12
$ws.Cells.Item(39,14).Value = 1.1388
$ws.Cells.Item(39,15).Value = 2.19
$ws.Cells.Item(39,16).Value = 0.929
$ws.Cells.Item(39,17).Value = $null
$ws.Cells.Item(39,18).Value = $null
$ws.Cells.Item(39,19).Value = 6
$ws.Cells.Item(39,20).Value = 1.2264
$ws.Cells.Item(39,21).Value = 1.1388
$ws.Cells.Item(39,22).Value = 1.0074
$ws.Cells.Item(39,23).Value = 0.8979
$ws.Cells.Item(39,24).Value = 0.7665
$ws.Cells.Item(39,25).Value = 0.6132
$ws.Cells.Item(39,26).Value = 0.4599
$ws.Cells.Item(39,27).Value = 0.4161
$ws.Cells.Item(39,28).Value = 1.155676
$ws.Cells.Item(39,29).Value = 'B2'
$ws.Cells.Item(39,30).Value = 1.1388
$ws.Cells.Item(39,31).Value = 1.12

# Row 40
$ws.Cells.Item(40,1).NumberFormat = "@"
$ws.Cells.Item(40,1).Value = '2025-10-25'
$ws.Cells.Item(40,1).ClearFormats()
$ws.Cells.Item(40,2).Value = 1.139
$ws.Cells.Item(40,3).Value = 1.22
$ws.Cells.Item(40,4).Value = 1.12
$ws.Cells.Item(40,5).Value = 1.131
$ws.Cells.Item(40,6).Value = 'wait'
$ws.Cells.Item(40,7).Value = $false
$ws.Cells.Item(40,8).Value = $null
$ws.Cells.Item(40,9).Value = 'SELL S'
$ws.Cells.Item(40,10).Value = 'HIGH'
$ws.Cells.Item(40,11).Value = $null
$ws.Cells.Item(40,12).Value = $null
$ws.Cells.Item(40,13).Value = 1.089717
$ws.Cells.Item(40,14).Value = 1.139
$ws.Cells.Item(40,15).Value = 2.19
$ws.Cells.Item(40,16).Value = 0.929
$ws.Cells.Item(40,17).Value = $null
$ws.Cells.Item(40,18).Value = 17.3
$ws.Cells.Item(40,19).Value = 6
$ws.Cells.Item(40,20).Value = 1.2264
$ws.Cells.Item(40,21).Value = 1.1388
$ws.Cells.Item(40,22).Value = 1.0074
$ws.Cells.Item(40,23).Value = 0.8979
$ws.Cells.Item(40,24).Value = 0.7665
$ws.Cells.Item(40,25).Value = 0.6132
$ws.Cells.Item(40,26).Value = 0.4599
$ws.Cells.Item(40,27).Value = 0.4161
$ws.Cells.Item(40,28).Value = 1.139
$ws.Cells.Item(40,29).Value = $null
$ws.Cells.Item(40,30).Value = $null
$ws.Cells.Item(40,31).Value = $null

# Row 41
$ws.Cells.Item(41,1).NumberFormat = "@"
$ws.Cells.Item(41,1).Value = '2025-10-25'
$ws.Cells.Item(41,1).ClearFormats()
$ws.Cells.Item(41,2).Value = 1.139
$ws.Cells.Item(41,3).Value = 1.22
$ws.Cells.Item(41,4).Value = 1.12
$ws.Cells.Item(41,5).Value = 1.131
$ws.Cells.Item(41,6).Value = 'wait'
$ws.Cells.Item(41,7).Value = $false
$ws.Cells.Item(41,8).Value = $null
$ws.Cells.Item(41,9).Value = $null
$ws.Cells.Item(41,10).Value = $null
$ws.Cells.Item(41,11).Value = $null
$ws.Cells.Item(41,12).Value = $null
$ws.Cells.Item(41,13).Value = $null
$ws.Cells.Item(41,14).Value = $null
$ws.Cells.Item(41,15).Value = 2.19
$ws.Cells.Item(41,16).Value = 0.929
$ws.Cells.Item(41,17).Value = 31.324004
$ws.Cells.Item(41,18).Value = 17.3
$ws.Cells.Item(41,19).Value = 6
$ws.Cells.Item(41,20).Value = 1.2264
$ws.Cells.Item(41,21).Value = 1.1388
$ws.Cells.Item(41,22).Value = 1.0074
$ws.Cells.Item(41,23).Value = 0.8979
$ws.Cells.Item(41,24).Value = 0.7665
$ws.Cells.Item(41,25).Value = 0.6132
$ws.Cells.Item(41,26).Value = 0.4599
$ws.Cells.Item(41,27).Value = 0.4161
$ws.Cells.Item(41,28).Value = 1.139
$ws.Cells.Item(41,29).Value = 'B2'
$ws.Cells.Item(41,30).Value = 1.1388
$ws.Cells.Item(41,31).Value = 1.12
